# Add an "Actual" results column (D) to the login test data sheet,
# matching formatting of the existing table, move the selection, and
# touch page setup so Excel regenerates it without the stale printer
# settings link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell D1: same look as A1:C1 (yellow fill + border) ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value2 = "Actual"

# --- D2: bordered cell with white (invisible) font, left blank ---
$ws.Range("D2").Borders.ColorIndex = 1
$ws.Range("D2").Borders.LineStyle = 1
$ws.Range("D2").Font.ThemeColor = 2

# --- D3:D5: bordered cells, left blank ---
$ws.Range("D3").Borders.ColorIndex = 1
$ws.Range("D3").Borders.LineStyle = 1
$ws.Range("D4").Borders.ColorIndex = 1
$ws.Range("D4").Borders.LineStyle = 1
$ws.Range("D5").Borders.ColorIndex = 1
$ws.Range("D5").Borders.LineStyle = 1

# --- keep the custom row heights, marking them as explicit/custom ---
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75

# --- move the active selection ---
$null = $ws.Range("F4").Select()

# --- touch page setup so the printer-settings relationship is dropped ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "done"
